$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price record was inserted as row 135 (Macroferia Regional de Talca -
# Brócoli, weekly data), pushing all subsequent rows (old 135-194) down by
# one to become rows 136-195.
$ws.Rows.Item(135).Insert()

$ws.Cells.Item(135,1).Value2  = 5
$ws.Cells.Item(135,2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(135,3).Value2  = "Maule"
$ws.Cells.Item(135,4).Value2  = 44466
$ws.Cells.Item(135,5).Value2  = 7
$ws.Cells.Item(135,6).Value2  = 100112023
$ws.Cells.Item(135,7).Value2  = "Brócoli"
$ws.Cells.Item(135,8).Value2  = "Sin especificar"
$ws.Cells.Item(135,9).Value2  = "Primera"
$ws.Cells.Item(135,10).Value2 = 3000
$ws.Cells.Item(135,11).Value2 = 600
$ws.Cells.Item(135,12).Value2 = 600
$ws.Cells.Item(135,13).Value2 = 600
$ws.Cells.Item(135,14).Value2 = "$/unidad"
$ws.Cells.Item(135,15).Value2 = "Región del Maule"
$ws.Cells.Item(135,16).Value2 = 600
$ws.Cells.Item(135,17).Value2 = 1
$ws.Cells.Item(135,18).Value2 = "Hortaliza"
